# "Generate Report for Handback"
#
# This script updates the localization-status workbook to reflect a
# handback event:
#   - Status text changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (on the Overview sheet's zh-cn/de-de
#     status columns, and on the Status column of the zh-cn and de-de
#     report sheets - they all shared the same text).
#   - The zh-cn sheet's "Latest Handback DateTime" is refreshed.
#   - The de-de sheet's "Latest Handback DateTime" is refreshed.
#   - The previous "handback not latest" error message is cleared from the
#     Error Detail column on both language sheets (now in sync).
#   - A few columns are resized (status columns wider to fit the longer
#     message, Error Detail columns narrower now that they are empty).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $newStatus
$ws2.Range("K2").Value = "2016-08-26 04:47:38"
$ws2.Range("P2").Value = ""

# --- de-de sheet ---
$ws3.Range("C2").Value = $newStatus
$ws3.Range("K2").Value = "2016-08-26 04:47:45"
$ws3.Range("P2").Value = ""

# --- Column width adjustments ---
# Overview: widen the zh-cn (E) and de-de (F) status columns
$ws1.Columns.Item(5).ColumnWidth = 29.166666666666668
$ws1.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn: widen Status column (C), narrow Error Detail column (P)
$ws2.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws2.Columns.Item(16).ColumnWidth = 12.833333333333334

# de-de: widen Status column (C), narrow Error Detail column (P)
$ws3.Columns.Item(3).ColumnWidth = 29.166666666666668
$ws3.Columns.Item(16).ColumnWidth = 12.833333333333334
